$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "52.801.23"
$ws.Range("E2").Value = "  +1.19%  "
$ws.Range("D3").Value = "2.988.94"
$ws.Range("E3").Value = "  +3.16%  "
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").Value = "'358.53"
$ws.Range("E5").Value = "  +1.72%  "
$ws.Range("D6").Value = "'110.57"
$ws.Range("E6").Value = "  -2.05%  "
$ws.Range("E7").Value = "  +3.25%  "
$ws.Range("E9").Value = "  +2.30%  "
$ws.Range("D10").Value = "'39.52"
$ws.Range("E10").Value = "  -1.05%  "
$ws.Range("D11").Value = "'0.0883"
$ws.Range("E11").Value = "  +2.61%  "
$ws.Range("E12").Value = "  +1.65%  "
$ws.Range("D13").Value = "'19.63"
$ws.Range("E13").Value = "  -0.98%  "
$ws.Range("B14").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C14").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D14").Value = "3.456.18"
$ws.Range("E14").Value = "  +3.08%  "
$ws.Range("B15").Value = "Polkadot"
$ws.Range("C15").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D15").Value = "'7.85"
$ws.Range("E15").Value = "  +1.45%  "
$ws.Range("D16").Value = "2.985.01"
$ws.Range("E16").Value = "  +2.88%  "
$ws.Range("D17").Value = "'1.00"
$ws.Range("E17").Value = "  +0.74%  "
$ws.Range("D18").Value = "52.797.59"
$ws.Range("E18").Value = "  +1.15%  "
$ws.Range("E19").Value = "  +6.68%  "
$ws.Range("D20").Value = "'7.70"
$ws.Range("E20").Value = "  +0.80%  "
$ws.Range("D21").Value = "'13.99"
$ws.Range("E21").Value = "  -1.17%  "
$ws.Range("E22").Value = "  +1.74%  "
$ws.Range("D23").Value = "'273.16"
$ws.Range("E23").Value = "  +1.67%  "
$ws.Range("D24").Value = "'70.81"
$ws.Range("E24").Value = "  +0.00%  "
$ws.Range("D25").Value = "'2.82"
$ws.Range("E25").Value = "  +1.36%  "
$ws.Range("E26").Value = "  +4.23%  "
$ws.Range("B27").Value = "EthereumClassic"
$ws.Range("C27").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D27").Value = "'27.53"
$ws.Range("E27").Value = "  +3.12%  "
$ws.Range("B28").Value = "Filecoin"
$ws.Range("C28").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D28").Value = "'7.74"
$ws.Range("E28").Value = "  +16.90%  "
$ws.Range("E29").Value = "  +0.03%  "
$ws.Range("D30").Value = "'0.109"
$ws.Range("E30").Value = "  +6.33%  "
$ws.Range("D31").Value = "'10.57"
$ws.Range("E31").Value = "  -0.25%  "
$ws.Range("D32").Value = "'38.25"
$ws.Range("E32").Value = "  +2.09%  "
$ws.Range("D33").Value = "'6.15"
$ws.Range("E33").Value = "  -0.66%  "
$ws.Range("D34").Value = "'2.18"
$ws.Range("E34").Value = "  +3.71%  "
$ws.Range("D35").Value = "'52.49"
$ws.Range("E35").Value = "  -1.33%  "
$ws.Range("E36").Value = "  +0.20%  "
$ws.Range("D37").Value = "'0.999"
$ws.Range("E37").Value = "  -0.01%  "
$ws.Range("D38").Value = "'3.27"
$ws.Range("E38").Value = "  -1.57%  "
$ws.Range("E39").Value = "  +0.24%  "
$ws.Range("D40").Value = "'18.32"
$ws.Range("E40").Value = "  -2.71%  "
$ws.Range("D41").Value = "'2.79"
$ws.Range("E41").Value = "  +1.18%  "
$ws.Range("E42").Value = "  +3.15%  "
$ws.Range("D43").Value = "'23.72"
$ws.Range("E43").Value = "  +3.58%  "
$ws.Range("D44").Value = "'119.96"
$ws.Range("E44").Value = "  +0.16%  "
$ws.Range("E45").Value = "  -1.29%  "
$ws.Range("D46").Value = "'3.49"
$ws.Range("E46").Value = "  -0.51%  "
$ws.Range("D47").Value = "2.154.58"
$ws.Range("E47").Value = "  -0.90%  "
$ws.Range("E48").Value = "  -5.10%  "
$ws.Range("B49").Value = "TheGraph"
$ws.Range("C49").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D49").Value = "'0.248"
$ws.Range("E49").Value = "  -6.07%  "
$ws.Range("B50").Value = "BEAM"
$ws.Range("C50").Value = "https://coinranking.com/coin/cYYMfXF4u+beam-beam"
$ws.Range("D50").Value = "'0.0356"
$ws.Range("E50").Value = "  +2.63%  "
$ws.Range("D51").Value = "'0.936"
$ws.Range("E51").Value = "  -1.44%  "
